# Add new FingerTips indicators to the FT_indicators table (Table1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FT_indicators")

# Extend the table (Table1 / ListObject) by two rows so the table range and
# autofilter pick up the new data (A1:C14 -> A1:C16).
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Copy the number formatting used by the existing data rows (A7:B13, style
# index 1) onto the two new rows before filling in the values.
$ws.Range("A7:B7").Copy() | Out-Null
$ws.Range("A15:B15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New indicator rows.
$ws.Cells.Item(15, 1).Value = 34
$ws.Cells.Item(15, 2).Value = 92781
$ws.Cells.Item(15, 3).Value = "GP"

$ws.Cells.Item(16, 1).Value = 55
$ws.Cells.Item(16, 2).Value = 93605
$ws.Cells.Item(16, 3).Value = "LA"

# Match the author's final navigation state: FT_indicators becomes the
# active/selected sheet tab, with different last-used selections on each
# sheet.
$wsCover = $wb.Worksheets.Item("Cover")
$wsCover.Range("B24").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("G12").Select() | Out-Null
